$wb = $excel.ActiveWorkbook

# The workbook has two sheets that carry a duplicate copy of the same
# event rows: "展览" (sheet 1) and "全部类型" (sheet 4). Both need the
# "想去人数" (column F) counts bumped by 1 for the three affected events.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 64
    $ws.Range("F6").Value = 6947
    $ws.Range("F10").Value = 411
}
